$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bug #3 (new row 5): "Attempt to run Countdown timer..." --------------
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Attempt to run Countdown timer with any time cap"
$ws.Range("C5").Value = "Timer should countdown from cap to 0:00"
$ws.Range("E5").Value = "Rsoderberg"

# Row 5 should look like the other bug rows (wrap text, 2-line tall).
$ws.Rows(5).RowHeight = 28.8

# F5 starts out "unresolved" (red) - same look bug #2 had before it got
# fixed, so just clone F4's current (red) format onto it.
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null

# --- Bug #2 (row 4): tidy up text, now resolved ----------------------------
$ws.Range("B4").Value = "User spinner to choose between time cap selections"

# Back to bug #3's last cell.
$ws.Range("D5").Value = "Timer does not countdown correctly, instead it simply ticks down one second and stops"

# F4 flips from "unresolved" (red) to "resolved" (green, default font).
$rng4 = $ws.Range("F4")
$rng4.NumberFormat = "@"
$rng4.WrapText = $true
$rng4.Interior.Color = 5296274
$rng4.Borders.LineStyle = 1

# --- Selection --------------------------------------------------------------
$ws.Range("D5").Select() | Out-Null
